# Generate Report for Handoff
# Updates status text from "In Translation" to "Ready for handoff" and
# refreshes the related timestamps on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns (E2, F2) and the
# "Latest HO Xliff Generate Date" column (G2).
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-22 10:39:19"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2).
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-22 10:39:14"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2).
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-22 10:39:19"

# Widen the status columns so the longer text fits, matching the
# autofit-driven column width change in the target workbook. Excel quantizes
# ColumnWidth to whole-pixel steps, so 16.33 is the input that lands on the
# closest achievable stored width to the target 17.2159881591797.
$wsOverview.Range("E1:F1").ColumnWidth = 16.33
$wsZhCn.Range("C1").ColumnWidth = 16.33
$wsDeDe.Range("C1").ColumnWidth = 16.33
